$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SNV")

# Insert two new columns before column D (shifts existing D:K quarterly data to F:M)
$ws.Columns("D:E").Insert()

# Copy number/date formatting from column F (the former column D) into the new D:E columns
$ws.Columns("F:F").Copy()
$ws.Columns("D:E").PasteSpecial(-4122)
$excel.CutCopyMode = $false


$ws.Cells.Item(7, 4).Value = 43465
$ws.Cells.Item(7, 5).Value = 43373
$ws.Cells.Item(8, 4).Value = 357400
$ws.Cells.Item(8, 5).Value = 343900
$ws.Cells.Item(9, 4).Value = "NA"
$ws.Cells.Item(9, 5).Value = "NA"
$ws.Cells.Item(10, 4).Value = "NA"
$ws.Cells.Item(10, 5).Value = "NA"
$ws.Cells.Item(12, 4).Value = "NA"
$ws.Cells.Item(12, 5).Value = "NA"
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(17, 4).Value = 71600
$ws.Cells.Item(17, 5).Value = 67300
$ws.Cells.Item(18, 4).Value = 285800
$ws.Cells.Item(18, 5).Value = 276600
$ws.Cells.Item(20, 4).Value = -141900
$ws.Cells.Item(20, 5).Value = -148600
$ws.Cells.Item(21, 4).Value = 157300
$ws.Cells.Item(21, 5).Value = 141100
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(23, 4).Value = 143900
$ws.Cells.Item(23, 5).Value = 128000
$ws.Cells.Item(24, 4).Value = 48600
$ws.Cells.Item(24, 5).Value = 18900
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(26, 4).Value = 95200
$ws.Cells.Item(26, 5).Value = 109100
$ws.Cells.Item(27, 4).Value = 92100
$ws.Cells.Item(27, 5).Value = 99300
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(29, 4).Value = 9900
$ws.Cells.Item(29, 5).Value = "NA"
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(32, 4).Value = 141900
$ws.Cells.Item(32, 5).Value = 148600
$ws.Cells.Item(33, 4).Value = 101900
$ws.Cells.Item(33, 5).Value = 99300
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(35, 4).Value = 101900
$ws.Cells.Item(35, 5).Value = 99300
$ws.Cells.Item(38, 4).Value = 43465
$ws.Cells.Item(38, 5).Value = 43373
$ws.Cells.Item(41, 4).Value = 468400
$ws.Cells.Item(41, 5).Value = 436500
$ws.Cells.Item(42, 4).Value = 675100
$ws.Cells.Item(42, 5).Value = 575400
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 0
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(48, 4).Value = 434300
$ws.Cells.Item(48, 5).Value = 431000
$ws.Cells.Item(49, 4).Value = 67200
$ws.Cells.Item(49, 5).Value = 67500
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 0
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(52, 4).Value = 141100
$ws.Cells.Item(52, 5).Value = 185100
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(54, 4).Value = 32669200
$ws.Cells.Item(54, 5).Value = 32075100
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 5).Value = 0
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(59, 4).Value = 0
$ws.Cells.Item(59, 5).Value = 0
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 0
$ws.Cells.Item(61, 4).Value = 1657200
$ws.Cells.Item(61, 5).Value = 1656900
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(66, 4).Value = 29535600
$ws.Cells.Item(66, 5).Value = 29035000
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(70, 4).Value = 195100
$ws.Cells.Item(70, 5).Value = 195100
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 0
$ws.Cells.Item(72, 4).Value = 843800
$ws.Cells.Item(72, 5).Value = 770800
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 0
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 0
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 0
$ws.Cells.Item(76, 4).Value = 2938500
$ws.Cells.Item(76, 5).Value = 2844900
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 0
$ws.Cells.Item(80, 4).Value = 43465
$ws.Cells.Item(80, 5).Value = 43373
$ws.Cells.Item(81, 4).Value = 101900
$ws.Cells.Item(81, 5).Value = 99300
$ws.Cells.Item(83, 4).Value = 13500
$ws.Cells.Item(83, 5).Value = 13100
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 0
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 0
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 0
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 0
$ws.Cells.Item(89, 4).Value = 164900
$ws.Cells.Item(89, 5).Value = 139800
$ws.Cells.Item(91, 4).Value = -14100
$ws.Cells.Item(91, 5).Value = -12300
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(94, 4).Value = -464000
$ws.Cells.Item(94, 5).Value = -449500
$ws.Cells.Item(96, 4).Value = -29200
$ws.Cells.Item(96, 5).Value = -29500
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(100, 4).Value = 430700
$ws.Cells.Item(100, 5).Value = 229800
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(102, 4).Value = 131600
$ws.Cells.Item(102, 5).Value = -79900

# Disable concurrent/multi-threaded calculation to match workbook calc settings
$excel.MultiThreadedCalculation.Enabled = $false
